# Update the "Generate Report for Handback" timestamps on the three sheets.
# These cells hold plain-text timestamps (shared strings), so assigning the
# new text keeps them as text rather than converting to a numeric date.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-30 19:21:21"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-30 19:21:17"

# zh-cn!K2 - Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-30 19:21:34"

# de-de!H2 - Correspond Handoff Datetime (shares the same text as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-30 19:21:21"

# de-de!K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-30 19:21:41"
